$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 223, shifting the existing rows 223-339
# down to 224-340 (and carrying their formatting with them).
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the new weekly price entry.
$ws.Cells.Item(223, 1).Value = 10
$ws.Cells.Item(223, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(223, 3).Value = "La Araucanía"
$ws.Cells.Item(223, 4).Value = 44518
$ws.Cells.Item(223, 5).Value = 9
$ws.Cells.Item(223, 6).Value = 100112028
$ws.Cells.Item(223, 7).Value = "Sandia"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 1030
$ws.Cells.Item(223, 11).Value = 900
$ws.Cells.Item(223, 12).Value = 1000
$ws.Cells.Item(223, 13).Value = 963
$ws.Cells.Item(223, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(223, 15).Value = "Perú"
$ws.Cells.Item(223, 16).Value = 963
$ws.Cells.Item(223, 17).Value = 1
$ws.Cells.Item(223, 18).Value = "Hortaliza"
